# "Automation UI Fixes for 17.11.1_RC web"
#
# The "Transactions" tab had a stray re-computed row (ID 34, Accrual,
# 01-Apr-2015, 7.51) at the top of its data. Switch to that sheet, select
# the whole row and delete it - every row below shifts up one position
# (ids/dates/types/amounts all move with their row), the used range shrinks
# by one row, and the sheet ends up the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

# Make "Transactions" the active sheet/tab.
$ws.Activate()

# Select row 2 in full (mirrors right-clicking the row header) then delete
# it, shifting everything below up by one row.
[void]$ws.Rows(2).Select()
$ws.Rows(2).Delete()
